$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.5149
$ws.Range("B3").Value = 6.024600000000004
$ws.Range("B14").Value = 5.608000000000003
$ws.Range("B16").Value = 6.186200000000003
$ws.Range("C18").Value = -11.6166
$ws.Range("B21").Value = 9.402100000000008
$ws.Range("B23").Value = 9.092100000000006
$ws.Range("C24").Value = -13.14429999999999
$ws.Range("B25").Value = 6.246899999999997
$ws.Range("C25").Value = -12.8653
$ws.Range("B26").Value = 5.717100000000006
$ws.Range("C27").Value = -12.7688
$ws.Range("B29").Value = 5.069100000000001
$ws.Range("C30").Value = -13.39049999999999
$ws.Range("C31").Value = -13.6063
$ws.Range("C39").Value = -12.14690000000001
$ws.Range("B40").Value = 8.916200000000002
$ws.Range("C42").Value = -12.45069999999999
$ws.Range("C48").Value = -11.26649999999999
$ws.Range("C51").Value = -11.1486
$ws.Range("C52").Value = -11.28469999999999
$ws.Range("B53").Value = 5.276099999999998
$ws.Range("C55").Value = -13.7732
$ws.Range("C56").Value = -12.83219999999999
$ws.Range("B57").Value = 4.881299999999996
$ws.Range("C57").Value = -13.56709999999999
$ws.Range("B59").Value = 4.698899999999996
$ws.Range("C60").Value = -13.59099999999999
$ws.Range("B65").Value = 5.873600000000002
$ws.Range("B69").Value = 5.621299999999995
$ws.Range("C73").Value = -12.19820000000001
$ws.Range("C74").Value = -12.76750000000001
$ws.Range("B79").Value = 8.935400000000007
$ws.Range("B83").Value = 5.5534
$ws.Range("C89").Value = -10.58130000000001
$ws.Range("C90").Value = -12.6527
$ws.Range("B91").Value = 5.023099999999999
$ws.Range("C92").Value = -10.846
$ws.Range("B93").Value = 5.881099999999996
$ws.Range("B100").Value = 5.1153
